$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing it to be stored as text so that
# numeric-looking strings (e.g. "584.68", "0.0000183", "8.10") are not
# silently coerced into numbers (which would lose formatting / precision).
function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# ---- Row 2 - Bitcoin ----
Set-TextValue "D2" "65.434.99"
$ws.Range("E2").Value = "  +4.53%  "

# ---- Row 3 - Ethereum ----
Set-TextValue "D3" "2.575.50"
$ws.Range("E3").Value = "  +4.79%  "

# ---- Row 4 - TetherUSD ----
$ws.Range("E4").Value = "  -0.02%  "

# ---- Row 5 - BNB ----
Set-TextValue "D5" "584.68"
$ws.Range("E5").Value = "  +2.26%  "

# ---- Row 6 - Solana ----
Set-TextValue "D6" "154.52"
$ws.Range("E6").Value = "  +5.20%  "

# ---- Row 7 - USDC ----
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  -0.10%  "

# ---- Row 8 - XRP ----
Set-TextValue "D8" "0.542"
$ws.Range("E8").Value = "  +1.85%  "

# ---- Row 9 - LidoStakedEther ----
Set-TextValue "D9" "2.579.15"
$ws.Range("E9").Value = "  +4.86%  "

# ---- Row 10 - Dogecoin ----
Set-TextValue "D10" "0.114"
$ws.Range("E10").Value = "  +2.88%  "

# ---- Row 11 - TRON ----
$ws.Range("E11").Value = "  -1.39%  "

# ---- Row 12 - Toncoin -> Cardano (rows 12 and 13 swap order) ----
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue "D12" "0.360"
$ws.Range("E12").Value = "  +2.39%  "

# ---- Row 13 - Cardano -> Toncoin ----
$ws.Range("B13").Value = "Toncoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D13" "5.29"
$ws.Range("E13").Value = "  +0.15%  "

# ---- Row 14 - Avalanche ----
Set-TextValue "D14" "29.46"
$ws.Range("E14").Value = "  +2.44%  "

# ---- Row 15 - ShibaInu ----
Set-TextValue "D15" "0.0000183"
$ws.Range("E15").Value = "  +4.31%  "

# ---- Row 16 - WrappedliquidstakedEther2.0 ----
Set-TextValue "D16" "3.041.67"
$ws.Range("E16").Value = "  +4.70%  "

# ---- Row 17 - WrappedBTC ----
Set-TextValue "D17" "65.312.66"
$ws.Range("E17").Value = "  +4.42%  "

# ---- Row 18 - WrappedEther ----
Set-TextValue "D18" "2.580.54"
$ws.Range("E18").Value = "  +4.85%  "

# ---- Row 19 - Uniswap ----
Set-TextValue "D19" "8.10"
$ws.Range("E19").Value = "  +2.03%  "

# ---- Row 20 - Chainlink ----
Set-TextValue "D20" "11.12"
$ws.Range("E20").Value = "  +2.49%  "

# ---- Row 21 - BitcoinCash ----
Set-TextValue "D21" "358.18"
$ws.Range("E21").Value = "  +10.29%  "

# ---- Row 22 - Polkadot ----
Set-TextValue "D22" "4.33"
$ws.Range("E22").Value = "  +4.78%  "

# ---- Row 23 - SuiNetwork ----
$ws.Range("E23").Value = "  +3.19%  "

# ---- Row 24 - Dai ----
$ws.Range("E24").Value = "  +0.09%  "

# ---- Row 25 - Aptos ----
Set-TextValue "D25" "10.11"
$ws.Range("E25").Value = "  +1.26%  "

# ---- Row 26 - Litecoin ----
Set-TextValue "D26" "66.46"
$ws.Range("E26").Value = "  +1.75%  "

# ---- Row 27 - Bittensor ----
Set-TextValue "D27" "638.19"
$ws.Range("E27").Value = "  -0.33%  "

# ---- Row 28 - PEPE ----
$ws.Range("E28").Value = "  +9.01%  "

# ---- Row 29 - WrappedeETH ----
Set-TextValue "D29" "2.679.49"
$ws.Range("E29").Value = "  +3.60%  "

# ---- Row 30 - Fetch.AI ----
$ws.Range("E30").Value = "  +5.37%  "

# ---- Row 31 - Binance-PegBSC-USD ----
Set-TextValue "D31" "0.999"
$ws.Range("E31").Value = "  -0.19%  "

# ---- Row 32 - InternetComputer(DFINITY) ----
Set-TextValue "D32" "8.19"
$ws.Range("E32").Value = "  +3.53%  "

# ---- Row 33 - PancakeSwap ----
Set-TextValue "D33" "1.89"
$ws.Range("E33").Value = "  +3.46%  "

# ---- Row 34 - Kaspa ----
$ws.Range("E34").Value = "  +5.47%  "

# ---- Row 35 - FirstDigitalUSD ----
Set-TextValue "D35" "0.999"
$ws.Range("E35").Value = "  +0.05%  "

# ---- Row 36 - ImmutableX ----
Set-TextValue "D36" "1.58"
$ws.Range("E36").Value = "  +4.98%  "

# ---- Row 37 - NEARProtocol ----
Set-TextValue "D37" "4.93"
$ws.Range("E37").Value = "  +4.31%  "

# ---- Row 38 - RenderToken ----
Set-TextValue "D38" "5.69"
$ws.Range("E38").Value = "  +7.93%  "

# ---- Row 39 - dogwifhat ----
Set-TextValue "D39" "2.89"
$ws.Range("E39").Value = "  +6.82%  "

# ---- Row 40 - EthereumClassic ----
Set-TextValue "D40" "19.27"
$ws.Range("E40").Value = "  +4.00%  "

# ---- Row 41 - Monero ----
$ws.Range("E41").Value = "  +2.71%  "

# ---- Row 42 - PolygonEcosystemToken ----
Set-TextValue "D42" "0.375"
$ws.Range("E42").Value = "  +2.22%  "

# ---- Row 43 - Stacks ----
$ws.Range("E43").Value = "  +5.62%  "

# ---- Row 44 - OKB ----
Set-TextValue "D44" "41.91"
$ws.Range("E44").Value = "  +0.44%  "

# ---- Row 45 - Aave ----
Set-TextValue "D45" "162.28"
$ws.Range("E45").Value = "  +6.59%  "

# ---- Row 46 - BabyDogeCoin -> USDe (rows 46 and 47 swap order) ----
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D46" "0.999"
$ws.Range("E46").Value = "  -0.04%  "

# ---- Row 47 - USDe -> BabyDogeCoin ----
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D47" "0.0₆0307"
$ws.Range("E47").Value = "  -1.10%  "

# ---- Row 48 - WhiteBITCoin ----
Set-TextValue "D48" "15.87"
$ws.Range("E48").Value = "  +3.72%  "

# ---- Row 49 - Filecoin ----
Set-TextValue "D49" "3.71"
$ws.Range("E49").Value = "  +4.09%  "

# ---- Row 50 - InjectiveProtocol ----
Set-TextValue "D50" "21.74"
$ws.Range("E50").Value = "  +7.56%  "

# ---- Row 51 - Mantle ----
Set-TextValue "D51" "0.634"
$ws.Range("E51").Value = "  +4.90%  "
